$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: volume number and report date range ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Crime statistics table updates (rows 14-30) ---
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = -50
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 18
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -5.263157894736
$ws.Range("M16").Value = -65.384615384615
$ws.Range("N16").Value = -89.411764705882
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -50
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 51
$ws.Range("K17").Value = -29.411764705882
$ws.Range("L17").Value = -10
$ws.Range("M17").Value = 125
$ws.Range("N17").Value = -20
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -35.294117647058
$ws.Range("I18").Value = 41
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = -2.380952380952
$ws.Range("L18").Value = -6.818181818181
$ws.Range("M18").Value = -8.888888888888
$ws.Range("N18").Value = -86.73139158576
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -43.396226415094
$ws.Range("I19").Value = 66
$ws.Range("J19").Value = 101
$ws.Range("K19").Value = -34.653465346534
$ws.Range("L19").Value = -43.589743589743
$ws.Range("M19").Value = -21.428571428571
$ws.Range("N19").Value = -26.666666666666
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 108.333333333333
$ws.Range("I20").Value = 58
$ws.Range("J20").Value = 43
$ws.Range("K20").Value = 34.883720930232
$ws.Range("L20").Value = 26.086956521739
$ws.Range("M20").Value = 81.25
$ws.Range("N20").Value = -92.909535452322
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -45.454545454545
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = -18.965517241379
$ws.Range("I21").Value = 224
$ws.Range("J21").Value = 260
$ws.Range("K21").Value = -13.846153846153
$ws.Range("L21").Value = -17.037037037037
$ws.Range("M21").Value = -3.030303030303
$ws.Range("N21").Value = -84.43363446838
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -55.555555555555
$ws.Range("L22").Value = -33.333333333333
$ws.Range("M22").Value = -20
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = -20
$ws.Range("L23").Value = -46.666666666666
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -17.391304347826
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -12.087912087912
$ws.Range("I24").Value = 165
$ws.Range("J24").Value = 191
$ws.Range("K24").Value = -13.612565445026
$ws.Range("L24").Value = -37.5
$ws.Range("M24").Value = 42.241379310344
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 63
$ws.Range("J25").Value = 86
$ws.Range("K25").Value = -26.744186046511
$ws.Range("L25").Value = -27.586206896551
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -78.571428571428
$ws.Range("F26").Value = 45
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 7.142857142857
$ws.Range("I26").Value = 97
$ws.Range("J26").Value = 73
$ws.Range("K26").Value = 32.876712328767
$ws.Range("L26").Value = 6.593406593406
$ws.Range("M26").Value = 59.016393442622
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = -57.142857142857
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 7
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = 57.142857142857
$ws.Range("L28").Value = 83.333333333333
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0

# --- Fix number formats / styles for cells that switched between text("N/A") and numeric ---
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("G28").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("H28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
